# Generate Report for Handoff
#
# The localization status flipped from "In Translation" to
# "Ready for handoff" and the handoff timestamps were refreshed. Widen the
# status columns so the longer text fits.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refreshed handoff / generate timestamps ---
$zhcn.Range("H2").Value = "2016-08-31 13:00:53"
$dede.Range("H2").Value = "2016-08-31 13:01:12"
$overview.Range("G2").Value = "2016-08-31 13:01:12"

# --- Widen the Status columns to fit "Ready for handoff" ---
$overview.Range("E1").EntireColumn.ColumnWidth = 16.33
$overview.Range("F1").EntireColumn.ColumnWidth = 16.33
$zhcn.Range("C1").EntireColumn.ColumnWidth = 16.33
$dede.Range("C1").EntireColumn.ColumnWidth = 16.33
